$wb = $excel.ActiveWorkbook

# --- Update Metadata sheet: Date value ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value = "2024-09-03T11:48:25+00:00"

# --- Update Elements sheet: Condition(s) / Mapping columns ---
$wsElem = $wb.Worksheets.Item("Elements")

# AI column = "Condition(s)" -> clear to empty string for rows 4, 6, 8, 9, 10
$wsElem.Range("AI4").Value = ""
$wsElem.Range("AI6").Value = ""
$wsElem.Range("AI8").Value = ""
$wsElem.Range("AI9").Value = ""
$wsElem.Range("AI10").Value = ""

# AK column = "Mapping: RIM Mapping" -> row 8 changes from "N/A" to "n/a"
$wsElem.Range("AK8").Value = "n/a"
